# Updated cryptos list on Tue May  2 19:37:40 UTC 2023 with GitHub Actions
#
# The "Price" (D) column stores numeric-looking values as TEXT in the
# source workbook (t="inlineStr"). Plain `.Value = "1.005"` assignment
# would be auto-coerced to a Number by Excel, which both changes the
# stored type and can silently drop significant trailing zeros
# (e.g. "0.9780" -> 0.978). A leading apostrophe forces Excel to keep
# the literal text entry (shows as "number stored as text"), which is
# used below for every Price cell whose new value parses as a number.
# Percent (E) cells already contain spaces/"%" so they stay text as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.888.30"
$ws.Range("E2").Value = "  +2.97%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.881.84"
$ws.Range("E3").Value = "  +3.15%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.18%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'327.42"
$ws.Range("E5").Value = "  -0.03%  "

# Row 6 - USDC
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  +0.13%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.4666"
$ws.Range("E7").Value = "  +1.19%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3956"
$ws.Range("E8").Value = "  +2.85%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.07931"
$ws.Range("E9").Value = "  +1.39%  "

# Row 10 - Polygon
$ws.Range("D10").Value = "'0.9780"
$ws.Range("E10").Value = "  +2.13%  "

# Row 11 - Solana
$ws.Range("D11").Value = "'22.44"
$ws.Range("E11").Value = "  +2.79%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.893.48"
$ws.Range("E12").Value = "  +8.61%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "'5.765"
$ws.Range("E13").Value = "  +1.71%  "

# Row 14 - Chainlink
$ws.Range("D14").Value = "'7.010"
$ws.Range("E14").Value = "  +2.13%  "

# Row 15 - TRON
$ws.Range("D15").Value = "'0.06994"
$ws.Range("E15").Value = "  +2.11%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "'88.92"
$ws.Range("E16").Value = "  +2.97%  "

# Row 17 - BinanceUSD
$ws.Range("D17").Value = "'1.005"
$ws.Range("E17").Value = "  +0.15%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "'0.00001015"
$ws.Range("E18").Value = "  +2.35%  "

# Row 19 - Avalanche
$ws.Range("D19").Value = "'17.02"
$ws.Range("E19").Value = "  +1.28%  "

# Row 20 - Dai (only Volume changes)
$ws.Range("E20").Value = "  +0.16%  "

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "28.879.29"
$ws.Range("E21").Value = "  +2.84%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'5.361"
$ws.Range("E22").Value = "  +0.86%  "

# Row 23 - Cosmos
$ws.Range("D23").Value = "'11.14"
$ws.Range("E23").Value = "  +1.62%  "

# Row 24 - Toncoin (only Volume changes)
$ws.Range("E24").Value = "  -0.47%  "

# Row 25 - WrappedliquidstakedEther2.0
$ws.Range("D25").Value = "2.073.46"
$ws.Range("E25").Value = "  +4.31%  "

# Row 26 - Monero
$ws.Range("D26").Value = "'153.57"
$ws.Range("E26").Value = "  +1.15%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'19.46"
$ws.Range("E27").Value = "  +1.76%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").Value = "'5.783"
$ws.Range("E28").Value = "  +1.43%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "'2.013"
$ws.Range("E29").Value = "  +2.10%  "

# Row 30 - BitcoinCash (only Volume changes)
$ws.Range("E30").Value = "  +3.04%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "'0.09399"
$ws.Range("E31").Value = "  +1.44%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "'0.9457"
$ws.Range("E32").Value = "  +0.93%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "'5.334"
$ws.Range("E33").Value = "  +1.47%  "

# Row 34 - ARBITRUM
$ws.Range("D34").Value = "'1.354"
$ws.Range("E34").Value = "  +3.59%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "'3.349"
$ws.Range("E35").Value = "  -2.26%  "

# Row 36 - Hedera
$ws.Range("D36").Value = "'0.05931"
$ws.Range("E36").Value = "  -0.67%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "'0.02126"
$ws.Range("E37").Value = "  -0.75%  "

# Row 38 - TrustWalletToken
$ws.Range("D38").Value = "'1.153"
$ws.Range("E38").Value = "  +0.61%  "

# Row 39 - FraxShare
$ws.Range("D39").Value = "'7.956"
$ws.Range("E39").Value = "  +5.36%  "

# Row 40 - TheSandbox (only Volume changes)
$ws.Range("E40").Value = "  +2.55%  "

# Row 41 - was Aptos, now Algorand
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.1799"
$ws.Range("E41").Value = "  +1.93%  "

# Row 42 - was Algorand, now Aptos
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'10.01"
$ws.Range("E42").Value = "  +1.02%  "

# Row 43 - Cronos
$ws.Range("D43").Value = "'0.07250"
$ws.Range("E43").Value = "  +3.66%  "

# Row 44 - EnergySwap
$ws.Range("D44").Value = "'11.85"
$ws.Range("E44").Value = "  +2.50%  "

# Row 45 - Decentraland
$ws.Range("D45").Value = "'0.5356"
$ws.Range("E45").Value = "  +2.07%  "

# Row 46 - was WEMIXToken, now RenderToken
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'2.145"
$ws.Range("E46").Value = "  -3.82%  "

# Row 47 - was RenderToken, now WEMIXToken
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'1.140"
$ws.Range("E47").Value = "  -8.46%  "

# Row 48 - NEARProtocol
$ws.Range("D48").Value = "'1.858"
$ws.Range("E48").Value = "  +1.78%  "

# Row 49 - Quant
$ws.Range("D49").Value = "'114.22"
$ws.Range("E49").Value = "  +1.75%  "

# Row 50 - MXToken (only Volume changes)
$ws.Range("E50").Value = "  +2.52%  "

# Row 51 - EOS
$ws.Range("D51").Value = "'1.033"
$ws.Range("E51").Value = "  +2.62%  "
